$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "id number" column is inserted right before "position" (old column F).
# Every header from the old "position" column onward shifts one column to the
# right (F->G, G->H, ... M->N), and a brand new "note" header lands in N1
# (phone moves from L to M... etc). Row 2's pre-formatted (but empty) cells in
# H2/I2 are left exactly where they are - this was not a structural
# "insert column" edit, just new/retyped header cells - so we set each header
# cell's value/format explicitly instead of using Columns.Insert().

$ws.Range("F1").Value2 = "id number"
$ws.Range("F1").NumberFormat = "General"

$ws.Range("G1").Value2 = "position"
$ws.Range("G1").NumberFormat = "General"

$ws.Range("H1").Value2 = "department"
$ws.Range("H1").NumberFormat = "General"

$ws.Range("I1").Value2 = "company"
$ws.Range("I1").NumberFormat = "yyyy\-mm\-dd;@"

$ws.Range("J1").Value2 = "joined date"
$ws.Range("J1").NumberFormat = "mm/dd/yy"

$ws.Range("K1").Value2 = "vaccineFirstDate"
$ws.Range("K1").NumberFormat = "mm/dd/yy"

$ws.Range("L1").Value2 = "vaccineSecondDate"
$ws.Range("L1").NumberFormat = "mm/dd/yy"

$ws.Range("M1").Value2 = "phone"
$ws.Range("M1").NumberFormat = "General"

$ws.Range("N1").Value2 = "note"
$ws.Range("N1").NumberFormat = "General"

# Resize the affected columns (F..N) to match the new, slightly narrower
# layout that makes room for the extra "id number" column.
# ColumnWidth is expressed in "characters"; the saved xlsx <col width="..">
# is ColumnWidth + 5/6, so subtract that offset to hit the target widths.
$ws.Columns.Item(6).ColumnWidth  = 9.37  - (5/6)
$ws.Columns.Item(7).ColumnWidth  = 10.69 - (5/6)
$ws.Columns.Item(8).ColumnWidth  = 11.46 - (5/6)
$ws.Columns.Item(9).ColumnWidth  = 10.25 - (5/6)
$ws.Columns.Item(10).ColumnWidth = 10.36 - (5/6)
$ws.Columns.Item(11).ColumnWidth = 14.22 - (5/6)
$ws.Columns.Item(12).ColumnWidth = 17.31 - (5/6)
$ws.Columns.Item(13).ColumnWidth = 14.22 - (5/6)
$ws.Columns.Item(14).ColumnWidth = 12.57 - (5/6)
